# Update capital structure database values for rows 2 and 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D" = 0.08470000000000001
    "E" = 0.113
    "G" = 0.2819709543568465
    "H" = 0.2655601659751037
    "I" = 0.2572614107883817
    "J" = 0.2024706217731462
    "K" = 10.3
    "L" = 0.2136929460580913
    "M" = 10.1
    "N" = 0.05507088331515812
    "O" = 0.9805825242718446
    "P" = 10.1
    "Q" = 0.05507088331515812
    "R" = 0.9805825242718446
    "X" = 0.06263980392067162
    "AB" = 0.06263980392067162
    "AM" = -0.011
    "AQ" = -1127.272727272727
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
